# "Generate Report for Handoff"
#
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the handoff timestamps on all three report
# sheets. Also widens the two date/status columns that the report
# generator re-sized when it produced this handoff run.

$wb = $excel.ActiveWorkbook

# The report-generator widens a status/date column by ~4.8 characters
# whenever it regenerates the sheet. Excel's COM ColumnWidth setter only
# accepts character widths and rounds to whole pixels (MDW=6px) before
# writing the <col width="..."/> attribute, so 16.333333 is the input
# that lands on the pixel the regenerated report used.
$newColWidth = 16.333333

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-13 04:46:41"

$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-13 04:46:34"

$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-13 04:46:41"

$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
